# Apply the two changes described by the commit:
#  1. Slide 6's table: switch its table style from the custom
#     "Table_0" style to the built-in "Medium Style 2 - Accent 1"
#     style (GUID {DF440677-253B-40FB-B9DC-BBC96329B62F}).
#  2. Re-colour the deck's main theme (ppt/theme/theme1.xml, the
#     theme used by the slide master / the whole deck) so that it
#     matches the stock Office "Office Theme" colour scheme instead
#     of the "Integral" colour scheme it previously had.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{DF440677-253B-40FB-B9DC-BBC96329B62F}")
    }
}

# --- 2. Theme colours -------------------------------------------------
# Office default theme colour scheme (dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink)
# expressed as OLE (BGR-packed) RGB integers, as required by
# ThemeColorScheme.Item(n).RGB.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
